# Applies the "cryptos list" data refresh described by the commit:
#   "Updated cryptos list on Sun Aug 18 20:54:42 UTC 2024 with GitHub Actions"
#
# For each affected row, Price (column D) and/or Volume(1h) (column E) are
# refreshed with the latest scraped values. Rows 47-49 additionally have their
# Coin/Link (columns B/C) reshuffled (VeChain, Maker, RenderToken rotate).
#
# All of these cells hold plain text (not numbers/dates) in the source sheet,
# e.g. D2 "59.725.65" is a literal string, not a date or a number. Assigning a
# numeric-looking string straight to .Value makes Excel coerce it into a real
# number/date, so those assignments are entered with a leading apostrophe
# (forces literal text, like typing it in the UI) and the cell style is put
# back to "Normal" afterwards, since the apostrophe entry also nudges Excel to
# stamp the cell with an explicit "Text" number format otherwise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.705.66"
$ws.Range("E2").Value = "  +0.63%  "
# Row 3
$ws.Range("D3").Value = "2.649.11"
$ws.Range("E3").Value = "  +1.57%  "
# Row 4
$ws.Range("E4").Value = "  +0.03%  "
# Row 5
$ws.Range("D5").Value = "'537.07"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = "'146.65"
$ws.Range("E6").Value = "  +3.92%  "
$ws.Range("D6").Style = "Normal"
# Row 7
$ws.Range("E7").Value = "  -0.04%  "
# Row 8
$ws.Range("E8").Value = "  +1.29%  "
# Row 9
$ws.Range("E9").Value = "  +5.95%  "
# Row 10
$ws.Range("D10").Value = "'0.102"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D10").Style = "Normal"
# Row 11
$ws.Range("E11").Value = "  +1.47%  "
# Row 12
$ws.Range("E12").Value = "  +0.07%  "
# Row 13
$ws.Range("D13").Value = "3.122.02"
$ws.Range("E13").Value = "  +1.82%  "
# Row 14
$ws.Range("D14").Value = "59.630.53"
$ws.Range("E14").Value = "  +0.63%  "
# Row 15
$ws.Range("D15").Value = "'21.42"
$ws.Range("E15").Value = "  +4.16%  "
$ws.Range("D15").Style = "Normal"
# Row 16
$ws.Range("D16").Value = "2.672.30"
$ws.Range("E16").Value = "  +2.08%  "
# Row 17
$ws.Range("E17").Value = "  +1.10%  "
# Row 18
$ws.Range("D18").Value = "'4.48"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("D18").Style = "Normal"
# Row 19
$ws.Range("D19").Value = "'340.39"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D19").Style = "Normal"
# Row 20
$ws.Range("D20").Value = "'10.35"
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("D20").Style = "Normal"
# Row 21
$ws.Range("D21").Value = "'6.22"
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("D21").Style = "Normal"
# Row 22
$ws.Range("E22").Value = "  +0.06%  "
# Row 23
$ws.Range("D23").Value = "'66.60"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D23").Style = "Normal"
# Row 24
$ws.Range("D24").Value = "'0.418"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("D24").Style = "Normal"
# Row 25
$ws.Range("E25").Value = "  -0.27%  "
# Row 26
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D26").Style = "Normal"
# Row 27
$ws.Range("D27").Value = "'7.32"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D27").Style = "Normal"
# Row 28
$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("E28").Value = "  +1.45%  "
# Row 29
$ws.Range("E29").Value = "  -0.08%  "
# Row 30
$ws.Range("E30").Value = "  -3.31%  "
# Row 31
$ws.Range("E31").Value = "  +1.67%  "
# Row 32
$ws.Range("D32").Value = "'18.90"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D32").Style = "Normal"
# Row 33
$ws.Range("D33").Value = "'150.61"
$ws.Range("D33").Style = "Normal"
# Row 34
$ws.Range("E34").Value = "  +1.00%  "
# Row 35
$ws.Range("E35").Value = "  +2.96%  "
# Row 36
$ws.Range("D36").Value = "'0.839"
$ws.Range("E36").Value = "  +3.09%  "
$ws.Range("D36").Style = "Normal"
# Row 37
$ws.Range("D37").Value = "'0.842"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("D37").Style = "Normal"
# Row 38
$ws.Range("E38").Value = "  -0.93%  "
# Row 39
$ws.Range("E39").Value = "  +1.66%  "
# Row 40
$ws.Range("D40").Value = "'285.85"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("D40").Style = "Normal"
# Row 41
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D41").Style = "Normal"
# Row 42
$ws.Range("E42").Value = "  +1.75%  "
# Row 43
$ws.Range("E43").Value = "  +0.14%  "
# Row 44
$ws.Range("D44").Value = "'0.0539"
$ws.Range("E44").Value = "  +2.65%  "
$ws.Range("D44").Style = "Normal"
# Row 45
$ws.Range("D45").Value = "'19.33"
$ws.Range("E45").Value = "  +3.74%  "
$ws.Range("D45").Style = "Normal"
# Row 46
$ws.Range("E46").Value = "  -0.69%  "
# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0227"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D47").Style = "Normal"
# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.967.67"
$ws.Range("E48").Value = "  +1.19%  "
# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.58"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D49").Style = "Normal"
# Row 50
$ws.Range("D50").Value = "'18.45"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D50").Style = "Normal"
# Row 51
$ws.Range("D51").Value = "'111.95"
$ws.Range("E51").Value = "  +0.76%  "
$ws.Range("D51").Style = "Normal"
